$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 999.63635
$ws.Range("I11").Value = 999.63635
$ws.Range("K11").Value = 999.63635
$ws.Range("M11").Value = -859.63635
$ws.Range("H55").Value = 387.5
$ws.Range("I55").Value = 133.66667
$ws.Range("J55").Value = 539.8
$ws.Range("K55").Value = 133.66667
$ws.Range("L55").Value = 539.8
$ws.Range("M55").Value = 80.33332999999999
$ws.Range("N55").Value = -967.8
$ws.Range("H76").Value = 5923
$ws.Range("I76").Value = 6162.875
$ws.Range("J76").Value = 4004
$ws.Range("K76").Value = 6162.875
$ws.Range("L76").Value = 4004
$ws.Range("M76").Value = -5847.875
$ws.Range("N76").Value = -4634
$ws.Range("H79").Value = 5923
$ws.Range("I79").Value = 6162.875
$ws.Range("J79").Value = 4004
$ws.Range("K79").Value = 6162.875
$ws.Range("L79").Value = 4004
$ws.Range("M79").Value = -5070.875
$ws.Range("N79").Value = -6188
$ws.Range("H137").Value = 1186.6923
$ws.Range("I137").Value = 1160.8
$ws.Range("K137").Value = 3482.4
$ws.Range("M137").Value = -932.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13235.632
$ws.Range("I32").Value = 14877.125
$ws.Range("K32").Value = 14877.125
$ws.Range("M32").Value = -14590.125
$ws.Range("H45").Value = 1463.5
$ws.Range("I45").Value = 1349.6666
$ws.Range("J45").Value = 1805
$ws.Range("K45").Value = 1349.6666
$ws.Range("L45").Value = 1805
$ws.Range("M45").Value = -972.6666
$ws.Range("N45").Value = -2559
$ws.Range("H74").Value = 1190.6923
$ws.Range("I74").Value = 831
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 831
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = 43
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1190.6923
$ws.Range("I77").Value = 831
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 4155
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = 213
$ws.Range("N77").Value = -18736
$ws.Range("H97").Value = 823.26666
$ws.Range("I97").Value = 736.2727
$ws.Range("J97").Value = 1062.5
$ws.Range("K97").Value = 736.2727
$ws.Range("L97").Value = 1062.5
$ws.Range("M97").Value = -240.2727
$ws.Range("N97").Value = -2054.5
$ws.Range("H122").Value = 4679.5713
$ws.Range("I122").Value = 4566.6665
$ws.Range("K122").Value = 13699.9995
$ws.Range("M122").Value = -11249.9995
$ws.Range("H131").Value = 57001
$ws.Range("J131").Value = 57001
$ws.Range("L131").Value = 57001
$ws.Range("N131").Value = -67081
$ws.Range("H132").Value = 4915.5386
$ws.Range("I132").Value = 6507.5454
$ws.Range("K132").Value = 19522.6362
$ws.Range("M132").Value = -16992.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 336.7931
$ws.Range("I22").Value = 196.21739
$ws.Range("K22").Value = 196.21739
$ws.Range("M22").Value = 153.78261
$ws.Range("H31").Value = 1763.258
$ws.Range("I31").Value = 1452.2142
$ws.Range("J31").Value = 4666.3335
$ws.Range("K31").Value = 1452.2142
$ws.Range("L31").Value = 4666.3335
$ws.Range("M31").Value = -1157.2142
$ws.Range("N31").Value = -5256.3335
$ws.Range("H34").Value = 1763.258
$ws.Range("I34").Value = 1452.2142
$ws.Range("J34").Value = 4666.3335
$ws.Range("K34").Value = 1452.2142
$ws.Range("L34").Value = 4666.3335
$ws.Range("M34").Value = -1250.2142
$ws.Range("N34").Value = -5070.3335
$ws.Range("H74").Value = 34000
$ws.Range("J74").Value = 34000
$ws.Range("L74").Value = 34000
$ws.Range("N74").Value = -35748
$ws.Range("H77").Value = 34000
$ws.Range("J77").Value = 34000
$ws.Range("L77").Value = 102000
$ws.Range("N77").Value = -110736
$ws.Range("H86").Value = 3726.3333
$ws.Range("I86").Value = 2543.2
$ws.Range("K86").Value = 2543.2
$ws.Range("M86").Value = -1420.2
$ws.Range("H89").Value = 3726.3333
$ws.Range("I89").Value = 2543.2
$ws.Range("K89").Value = 12716
$ws.Range("M89").Value = -7100
$ws.Range("H100").Value = 10000000
$ws.Range("J100").Value = 10000000
$ws.Range("L100").Value = 10000000
$ws.Range("N100").Value = -10002164
$ws.Range("H122").Value = 2205.1904
$ws.Range("I122").Value = 1912.0667
$ws.Range("J122").Value = 2938
$ws.Range("K122").Value = 5736.2001
$ws.Range("L122").Value = 8814
$ws.Range("M122").Value = -3286.2001
$ws.Range("N122").Value = -13714
$ws.Range("H134").Value = 2225.7297
$ws.Range("I134").Value = 1547.6538
$ws.Range("J134").Value = 3828.4546
$ws.Range("K134").Value = 4642.9614
$ws.Range("L134").Value = 11485.3638
$ws.Range("M134").Value = -2107.9614
$ws.Range("N134").Value = -16555.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 326.76923
$ws.Range("I5").Value = 326.76923
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 980.30769
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -868.30769
$ws.Range("N5").ClearContents()
$ws.Range("H23").Value = 229
$ws.Range("I23").Value = 83.75
$ws.Range("J23").Value = 270.5
$ws.Range("K23").Value = 251.25
$ws.Range("L23").Value = 811.5
$ws.Range("M23").Value = -16.25
$ws.Range("N23").Value = -1281.5
$ws.Range("H68").Value = 684.5833
$ws.Range("I68").Value = 504
$ws.Range("J68").Value = 744.7778
$ws.Range("K68").Value = 1512
$ws.Range("L68").Value = 2234.3334
$ws.Range("M68").Value = -701
$ws.Range("N68").Value = -3856.3334
$ws.Range("H71").Value = 684.5833
$ws.Range("I71").Value = 504
$ws.Range("J71").Value = 744.7778
$ws.Range("K71").Value = 4536
$ws.Range("L71").Value = 6703.000199999999
$ws.Range("M71").Value = -480
$ws.Range("N71").Value = -14815.0002
$ws.Range("H105").Value = 9942.857
$ws.Range("J105").Value = 9942.857
$ws.Range("L105").Value = 29828.571
$ws.Range("N105").Value = -35070.571
$ws.Range("H109").Value = 2618.3333
$ws.Range("I109").Value = 700
$ws.Range("J109").Value = 3577.5
$ws.Range("K109").Value = 2100
$ws.Range("L109").Value = 10732.5
$ws.Range("M109").Value = -1060
$ws.Range("N109").Value = -12812.5
$ws.Range("H119").Value = 2018.2858
$ws.Range("I119").Value = 2018.2858
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 6054.857400000001
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -1216.857400000001
$ws.Range("N119").ClearContents()
$ws.Range("H120").Value = 15385.429
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H121").Value = 1238.8334
$ws.Range("I121").Value = 600
$ws.Range("J121").Value = 1366.6
$ws.Range("K121").Value = 1800
$ws.Range("L121").Value = 4099.799999999999
$ws.Range("M121").Value = -490
$ws.Range("N121").Value = -6719.799999999999
$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 1248.5
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 11236.5
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -16136.5
$ws.Range("H125").Value = 3980
$ws.Range("J125").Value = 3980
$ws.Range("L125").Value = 11940
$ws.Range("N125").Value = -21780
$ws.Range("H132").Value = 1744.5294
$ws.Range("I132").Value = 1390.8
$ws.Range("K132").Value = 12517.2
$ws.Range("M132").Value = -9987.199999999999
$ws.Range("H135").Value = 326.76923
$ws.Range("I135").Value = 326.76923
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2940.92307
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -405.9230699999998
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2662.318
$ws.Range("I102").Value = 2245.3845
$ws.Range("J102").Value = 3264.5557
$ws.Range("K102").Value = 2245.3845
$ws.Range("L102").Value = 3264.5557
$ws.Range("M102").Value = -623.3845000000001
$ws.Range("N102").Value = -6508.5557
$ws.Range("H107").Value = 1428.5555
$ws.Range("I107").Value = 2270.8
$ws.Range("J107").Value = 375.75
$ws.Range("K107").Value = 2270.8
$ws.Range("L107").Value = 375.75
$ws.Range("M107").Value = -350.8000000000002
$ws.Range("N107").Value = -4215.75
$ws.Range("H122").Value = 4022.889
$ws.Range("J122").Value = 4150.75
$ws.Range("L122").Value = 12452.25
$ws.Range("N122").Value = -17352.25
$ws.Range("H131").Value = 48244.5
$ws.Range("J131").Value = 48244.5
$ws.Range("L131").Value = 48244.5
$ws.Range("N131").Value = -58324.5
$ws.Range("H132").Value = 2485.6875
$ws.Range("I132").Value = 1940.1111
$ws.Range("J132").Value = 3187.1428
$ws.Range("K132").Value = 5820.3333
$ws.Range("L132").Value = 9561.428400000001
$ws.Range("M132").Value = -3290.3333
$ws.Range("N132").Value = -14621.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1007.6842
$ws.Range("I16").Value = 1009.7647
$ws.Range("K16").Value = 1009.7647
$ws.Range("M16").Value = -839.7646999999999
$ws.Range("H122").Value = 8700609
$ws.Range("I122").Value = 3636.7273
$ws.Range("K122").Value = 10910.1819
$ws.Range("M122").Value = -8460.1819
$ws.Range("H136").Value = 34396452
$ws.Range("I136").Value = 45001840
$ws.Range("J136").Value = 1254625.6
$ws.Range("K136").Value = 135005520
$ws.Range("L136").Value = 3763876.8
$ws.Range("M136").Value = -135002970
$ws.Range("N136").Value = -3768976.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41669268
$ws.Range("I122").Value = 76925096
$ws.Range("J122").Value = 3292.9092
$ws.Range("K122").Value = 230775288
$ws.Range("L122").Value = 9878.7276
$ws.Range("M122").Value = -230772838
$ws.Range("N122").Value = -14778.7276
$ws.Range("H136").Value = 1863.1786
$ws.Range("I136").Value = 1746.4783
$ws.Range("J136").Value = 2400
$ws.Range("K136").Value = 5239.4349
$ws.Range("L136").Value = 7200
$ws.Range("M136").Value = -2689.4349
$ws.Range("N136").Value = -12300
